# Revert capacity charts to show kilowatts (instead of watts) on the y-axis.
# The underlying worksheet values are converted from watts to kilowatts
# (divide by 1000), the custom number format used for those cells gains a
# single decimal place, and the chart's value-axis title/number format are
# updated to reflect kilowatts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Worksheet data: convert the Wind (G) and Solar (E) capacity figures
#    from watts to kilowatts.
# ---------------------------------------------------------------------
$ws.Range("G3").Value = 18000 / 1000

$ws.Range("E13").Value = 10100 / 1000
$ws.Range("E14").Value = 482400 / 1000
$ws.Range("E15").Value = 64400 / 1000
$ws.Range("E16").Value = 10800 / 1000
$ws.Range("E17").Value = 17400 / 1000
$ws.Range("E18").Value = 21200 / 1000
$ws.Range("E19").Value = 45500 / 1000
$ws.Range("E20").Value = 14400 / 1000
$ws.Range("E21").Value = 47040 / 1000
$ws.Range("E22").Value = 81890 / 1000
$ws.Range("E23").Value = 114374 / 1000
$ws.Range("E24").Value = 258930 / 1000
$ws.Range("E25").Value = 252695 / 1000
$ws.Range("E26").Value = 1263990 / 1000

# ---------------------------------------------------------------------
# 2) Number format: the shared custom format (numFmtId 164, "#,##0") used
#    by all the capacity columns (B:G, rows 2-26) now shows one decimal
#    place.
# ---------------------------------------------------------------------
$ws.Range("B2:G26").NumberFormat = "#,##0.0"

# ---------------------------------------------------------------------
# 3) Chart: the value axis now reads "Kilowatts (kW)" with a plain
#    thousands-separated number format instead of the old "[K]" style.
# ---------------------------------------------------------------------
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$valueAxis = $chart.Axes(2)
$valueAxis.AxisTitle.Text = "Kilowatts (kW)"
$valueAxis.TickLabels.NumberFormat = "#,##0"
